$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, pushing existing rows 23-47 down to 24-48.
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with its data (same fixed market/category info as the
# rest of the sheet, with its own date/volume/price figures).
$ws.Cells.Item(23, 1).Value = 9
$ws.Cells.Item(23, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44757
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(24, 4).NumberFormat
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 100112035
$ws.Cells.Item(23, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 34
$ws.Cells.Item(23, 11).Value = 17000
$ws.Cells.Item(23, 12).Value = 18000
$ws.Cells.Item(23, 13).Value = 17500
$ws.Cells.Item(23, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(23, 15).Value = "Hijuelas"
$ws.Cells.Item(23, 16).Value = 1167
$ws.Cells.Item(23, 17).Value = 15
$ws.Cells.Item(23, 18).Value = "Hortaliza"
